{"js": "// Convert the complex field \"{ m:self.name }\" (fldChar begin / instrText /\n// fldChar end runs) into plain literal-text runs:\n//   {   m   :   self(orange)   .name}\n// while preserving the host paragraph's own identity and the orange color\n// that was applied to the \"self\" portion of the field code.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find the paragraph that owns the Word field (mirrors the template's single\n// \"{ m:self.name }\" field paragraph).\nfor (const p of paragraphs.items) {\n  p.fields.load(\"items\");\n}\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.fields.items.length > 0) {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  const ooxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n          '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' +\n              '<w:p w:rsidP=\"00F5495F\" w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\">' +\n                '<w:r><w:t>{</w:t></w:r>' +\n                '<w:r><w:t>m</w:t></w:r>' +\n                '<w:r><w:t>:</w:t></w:r>' +\n                '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>self</w:t></w:r>' +\n                '<w:r><w:t xml:space=\"preserve\">.name}</w:t></w:r>' +\n              '</w:p>' +\n            '</w:body>' +\n          '</w:document>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n    '</pkg:package>';\n\n  target.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Convert the complex field \"{ m:self.name }\" (fldChar begin/instrText/fldChar\n# end runs) into plain literal-text runs:  {  m  :  self(orange)  .name}\n# while preserving the host paragraph's own attributes and the orange color\n# that was applied to the \"self\" portion of the field code.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the Word field.\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Fields.Count -gt 0) {\n        $targetPara = $p\n        break\n    }\n}\nif ($targetPara -eq $null) {\n    $targetPara = $d.Fields.Item(1).Code.Paragraphs.Item(1)\n}\n\n$r = $targetPara.Range\n\n$xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"' `\n    + ' w:rsidP=\"00F5495F\" w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\">' `\n    + '<w:r><w:t>{</w:t></w:r>' `\n    + '<w:r><w:t>m</w:t></w:r>' `\n    + '<w:r><w:t>:</w:t></w:r>' `\n    + '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>self</w:t></w:r>' `\n    + '<w:r><w:t xml:space=\"preserve\">.name}</w:t></w:r>' `\n    + '</w:p>'\n\n$r.InsertXML($xml)\n"}
